# Add EDAM_DEF column (F) to Sheet1, matching the AFO_EDAM mapping update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same bold/bordered/centered style as the other header cells (copy format from E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "EDAM_DEF"

# EDAM_DEF values for rows 2-51 (row i -> $defs[i-2]), derived from the EDAM_IRI definition lookup.
$defs = @(
  '[''Assign molecular sequences, structures or other biological data to a specific group or category according to qualities it shares with that group or category.'']',
  '[''A text token, number or something else which identifies an entity, but which may not be persistent (stable) or unique (the same identifier may identify multiple things).'']',
  '[''One or more molecular sequences, possibly with associated annotation.'']',
  '[''The acquisition of data, typically measurements of physical systems using any type of sampling system, or by another other means.'']',
  '[''Compare two or more things to identify similarities.'']',
  '[''The visual representation of an object.'']',
  '[''An analytical technique that exploits the magenetic properties of certain atomic nuclei to provide information on the structure, dynamics, reaction state and chemical environment of molecules.'']',
  '[''Mathematical determination of the value of something, typically a properly of a  molecule.'']',
  '[''An analytical chemistry technique that measures the mass-to-charge ratio and abundance of ions in the gas phase.'']',
  '[''Apply analytical methods to existing data of a specific type.'']',
  '[''Apply analytical methods to existing data of a specific type.'']',
  '[''Construct some entity (typically a molecule sequence) from component pieces.'']',
  '[''Counting and measuring experimentally determined observations into quantities.'']',
  '[''Experimental methods for biomolecular structure determination, such as X-ray crystallography, nuclear magnetic resonance (NMR), circular dichroism (CD) spectroscopy, microscopy etc., including the assignment or modelling of molecular structure from such data.'']',
  '[''Sort a set of files or data items according to some property.'']',
  '[''The analysis of a image (typically a digital image) of some type in order to extract information from it.'']',
  '[''Generate an index of (typically a file of) biological data.'']',
  '[''Identify peaks in a spectrum from a mass spectrometry, NMR, or some other spectrum-generating experiment.'']',
  '[''Analyse data in order to deduce properties of an underlying distribution or population.'']',
  '[''One or more molecular sequences, possibly with associated annotation.'']',
  '[''A human-readable collection of information including annotation on a biological entity or phenomena, computer-generated reports of analysis of primary data (e.g. sequence or structural), and metadata (data about primary data) or any other free (essentially unformatted) text, as distinct from the primary data itself.'']',
  '[''Data (typically biological or biomedical) that has been rendered into an image, typically for display on screen.'']',
  '["Analyse DNA sequence data to identify differences between the genetic composition (genotype) of an individual compared to other individual''s or a reference sequence."]',
  '[''Convert a data set from one form to another.'']',
  '[''3D coordinate and associated data for a macromolecular tertiary (3D) structure or part of a structure.'']',
  '[''A function that processes a set of inputs and results in a set of outputs, or associates arguments (inputs) with values (outputs).'']',
  '[''Protein secondary structure from protein coordinate or circular dichroism (CD) spectroscopic data.'']',
  '[''The spectrum of frequencies of electromagnetic radiation emitted from a molecule as a result of some spectroscopy experiment.'']',
  '[''Spectra from mass spectrometry.'']',
  '[''One or more fingerprints (sequence classifiers) as used in the PRINTS database.'']',
  '[''An array of numerical values.'']',
  '[''Typically a simple numerical or string value that controls the operation of a tool.'']',
  '[''The name of a family of organism.'']',
  '[''A Uniform Resource Locator (URL).'']',
  '[''Name or other identifier of an entry from a biosample database.'']',
  '[''A string of characters that name or otherwise identify a resource on the Internet.'']',
  '[''The concentration of a chemical compound.'']',
  '[''The name (or part of a name) of a file (of any type).'']',
  '[''Typically a simple numerical or string value that controls the operation of a tool.'']',
  '[''A specification of a chemical structure.'']',
  '[''Raw data from or annotation on laboratory experiments.'']',
  '[''A digital data archive typically based around a relational model but sometimes using an object-oriented, tree or graph-based model.'']',
  '[''Mass of a molecule.'']',
  '[''Format used for images and image metadata.'']',
  '[''The extinction coefficient of a protein.'']',
  '[''Data on an error generated by computer system or tool.'']',
  '[''Validate some data.'']',
  '[''An array of numerical values.'']',
  '[''A valid email address of an end-user.'']',
  '[''Specification of one or more colors.'']'
)

for ($i = 0; $i -lt $defs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $defs[$i]
}
